$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 43 - fill in DPN (H), MFR (I), MPN (J) for the 10k resistor entry
$ws.Range("H43").Value = "311-10KGRCT-ND"
$ws.Range("I43").Value = "Yageo"
$ws.Range("J43").Value = "RC0603JR-0710KL"

# Row 49 - same values for the other 10k resistor entry
$ws.Range("H49").Value = "311-10KGRCT-ND"
$ws.Range("I49").Value = "Yageo"
$ws.Range("J49").Value = "RC0603JR-0710KL"

# Update view state: scroll/freeze pane position and active selection
$ws.Activate()
$ws.Range("A28").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G49:J49").Select()
